$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Strip the stray leading non-breaking-space character that precedes every
# country name in column A (rows 2-206), e.g. "\u00A0Great Britain" -> "Great Britain".
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

$rng = $ws.Range("A2:A" + $lastRow)
foreach ($cell in $rng.Cells) {
    $val = $cell.Value2
    if ($val -ne $null) {
        $trimmed = $val.TrimStart([char]0x00A0, ' ')
        if ($trimmed -ne $val) {
            $cell.Value2 = $trimmed
        }
    }
}

# Move the active selection to E6, matching the edited workbook's view state.
$ws.Range("E6").Select()
